# Fruta / hortaliza, semanal
# The commit reshuffles the (Fecha, Volumen, Precio minimo, Precio maximo,
# Precio promedio ponderado, Precio $/Kg) tuple across data rows 2..41 -
# i.e. row N ends up holding the values that used to live on row Map[N].
# All other columns (A,B,C,E,F,G,H,I,J,K,L,Q,R,T) stay put.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# target row -> source row (source row's ORIGINAL D/M/N/O/P/S values move to target row)
$map = @{
    2  = 19
    3  = 7
    4  = 8
    5  = 27
    6  = 29
    7  = 17
    8  = 6
    9  = 22
    10 = 34
    11 = 37
    12 = 5
    13 = 10
    14 = 14
    15 = 31
    16 = 16
    17 = 9
    18 = 35
    19 = 28
    20 = 20
    21 = 39
    22 = 30
    23 = 26
    24 = 33
    25 = 3
    26 = 24
    27 = 18
    28 = 21
    29 = 32
    30 = 12
    31 = 25
    32 = 23
    33 = 41
    34 = 2
    35 = 13
    36 = 4
    37 = 38
    38 = 36
    39 = 40
    40 = 11
    41 = 15
}

# Snapshot the original values for the columns that move, keyed by row,
# before any writes happen (so later writes never clobber a value we
# still need to read for another row).
$cols = @("D", "M", "N", "O", "P", "S")
$snapshot = @{}
foreach ($row in 2..41) {
    $vals = @{}
    foreach ($col in $cols) {
        $vals[$col] = $ws.Range("$col$row").Value()
    }
    $snapshot[$row] = $vals
}

foreach ($row in 2..41) {
    $srcRow = $map[$row]
    $vals = $snapshot[$srcRow]
    foreach ($col in $cols) {
        $ws.Range("$col$row").Value = $vals[$col]
    }
}
